$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.287.80'
$ws.Range("E2").Value = '  +0.37%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.329.36'
$ws.Range("E3").Value = '  -0.03%  '

# Row 4
$ws.Range("E4").Value = '  -0.17%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '188.89'
$ws.Range("E5").Value = '  +3.99%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '560.51'
$ws.Range("E6").Value = '  +0.37%  '

# Row 7
$ws.Range("E7").Value = '  -0.05%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.590'
$ws.Range("E8").Value = '  +0.33%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.321.30'
$ws.Range("E9").Value = '  -0.10%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.185'
$ws.Range("E10").Value = '  +0.53%  '

# Row 11
$ws.Range("E11").Value = '  +1.24%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '47.86'
$ws.Range("E12").Value = '  +1.50%  '

# Row 13
$ws.Range("E13").Value = '  +3.22%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.71'
$ws.Range("E14").Value = '  +2.15%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.857.76'
$ws.Range("E15").Value = '  -0.29%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '607.17'
$ws.Range("E16").Value = '  +1.06%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.284.46'
$ws.Range("E17").Value = '  +0.10%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.09'
$ws.Range("E18").Value = '  +0.98%  '

# Row 19
$ws.Range("E19").Value = '  +1.53%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.340.25'
$ws.Range("E20").Value = '  +0.10%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.18'
$ws.Range("E21").Value = '  -1.48%  '

# Row 22
$ws.Range("E22").Value = '  +1.98%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '18.52'
$ws.Range("E23").Value = '  +10.14%  '

# Row 24
$ws.Range("E24").Value = '  +1.69%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '100.77'
$ws.Range("E25").Value = '  +0.65%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.00'
$ws.Range("E26").Value = '  -0.25%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.04'
$ws.Range("E27").Value = '  +0.59%  '

# Row 28
$ws.Range("E28").Value = '  +5.20%  '

# Row 29
$ws.Range("E29").Value = '  +5.14%  '

# Row 30
$ws.Range("E30").Value = '  +0.30%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '30.59'
$ws.Range("E31").Value = '  +0.03%  '

# Row 32
$ws.Range("E32").Value = '  +9.69%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.86'
$ws.Range("E33").Value = '  +4.27%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '584.10'
$ws.Range("E34").Value = '  +9.91%  '

# Row 35
$ws.Range("E35").Value = '  +1.84%  '

# Row 36
$ws.Range("E36").Value = '  +2.14%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.730.48'
$ws.Range("E37").Value = '  -1.20%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '57.29'
$ws.Range("E38").Value = '  -0.95%  '

# Row 39
$ws.Range("E39").Value = '  +0.17%  '

# Row 40 - CoreDAO
$ws.Range("B40").Value = 'CoreDAO'
$ws.Range("C40").Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.52'
$ws.Range("E40").Value = '  +17.22%  '

# Row 41 - PEPE
$ws.Range("B41").Value = 'PEPE'
$ws.Range("C41").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0₃0732'
$ws.Range("E41").Value = '  +3.10%  '

# Row 42 - Kaspa
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.132'
$ws.Range("E42").Value = '  +6.00%  '

# Row 43 - InjectiveProtocol
$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '34.19'
$ws.Range("E43").Value = '  +7.94%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.30'
$ws.Range("E44").Value = '  -4.37%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.70'
$ws.Range("E45").Value = '  +1.70%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.345'
$ws.Range("E46").Value = '  +2.06%  '

# Row 47 - ApeXProtocol
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.35'
$ws.Range("E47").Value = '  +2.86%  '

# Row 48 - VeChain
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0426'
$ws.Range("E48").Value = '  +3.50%  '

# Row 49
$ws.Range("E49").Value = '  +1.40%  '

# Row 50
$ws.Range("E50").Value = '  +0.43%  '

# Row 51
$ws.Range("E51").Value = '  -0.01%  '
